# "updated .Net Environments for android 16"
#
# Adds a new row to the Android ".NET Environment" table for a
# Samsung A15 running Android 16, inserted as the new first data row
# (row 4). All existing data rows shift down by one, exactly like a
# normal Excel "Insert Row" above the current first entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above row 4 (the current first data row),
# pushing the existing rows 4-23 down to rows 5-24.
$ws.Rows.Item(4).Insert() | Out-Null

# Fill in the new row 4 with the Android 16 entry (column E / "Mono"
# is intentionally left blank, matching the other "summary" rows).
$ws.Range("A4").Value = "Samsung"
$ws.Range("B4").Value = "A15"
$ws.Range("C4").Value = "16"
$ws.Range("D4").Value = "5.15.180 Thu Oct 23 2025 00:27:49 UTC"
$ws.Range("F4").Value = "9.0.9"
$ws.Range("G4").Value = "ARM 64"

# Match the author's saved cell selection.
$ws.Range("D4").Select() | Out-Null
